$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13 (shifts existing rows 13-59 down by one)
$ws.Rows.Item(13).Insert()

$ws.Range("B13").Value = "AENJ"
$ws.Range("C13").Value = "Att. Enjambem."
